$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.835.82'
$ws.Range("E2").Value = '  -0.35%  '
$ws.Range("D3").Value = '2.354.89'
$ws.Range("E3").Value = '  -2.15%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = "'317.36"
$ws.Range("E5").Value = '  -5.79%  '
$ws.Range("D6").Value = "'106.97"
$ws.Range("E6").Value = '  +4.65%  '
$ws.Range("E7").Value = '  -1.67%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").Value = "'0.619"
$ws.Range("E9").Value = '  -3.28%  '
$ws.Range("E10").Value = '  +0.21%  '
$ws.Range("D11").Value = "'0.0925"
$ws.Range("E11").Value = '  -1.07%  '
$ws.Range("D12").Value = "'8.54"
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("D13").Value = "'1.00"
$ws.Range("E13").Value = '  -4.29%  '
$ws.Range("E14").Value = '  +0.06%  '
$ws.Range("D15").Value = "'15.90"
$ws.Range("E15").Value = '  -6.39%  '
$ws.Range("D16").Value = '2.711.08'
$ws.Range("E16").Value = '  -2.01%  '
$ws.Range("D17").Value = '2.357.39'
$ws.Range("E17").Value = '  -1.92%  '
$ws.Range("D18").Value = '42.805.83'
$ws.Range("E18").Value = '  -0.38%  '
$ws.Range("E19").Value = '  -0.45%  '
$ws.Range("D20").Value = "'0.0000106"
$ws.Range("E20").Value = '  -1.31%  '
$ws.Range("D21").Value = "'76.01"
$ws.Range("E21").Value = '  -0.87%  '
$ws.Range("E22").Value = '  -7.69%  '
$ws.Range("D23").Value = "'266.63"
$ws.Range("E23").Value = '  -1.69%  '
$ws.Range("D24").Value = "'2.30"
$ws.Range("E24").Value = '  -3.71%  '
$ws.Range("D25").Value = "'9.38"
$ws.Range("E25").Value = '  -9.22%  '
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").Value = "'11.39"
$ws.Range("E27").Value = '  -3.75%  '
$ws.Range("D28").Value = "'23.36"
$ws.Range("E28").Value = '  -3.69%  '
$ws.Range("E29").Value = '  +2.18%  '
$ws.Range("D30").Value = "'36.67"
$ws.Range("E30").Value = '  +0.58%  '
$ws.Range("D31").Value = "'167.34"
$ws.Range("E31").Value = '  -3.83%  '
$ws.Range("D32").Value = "'0.0908"
$ws.Range("E32").Value = '  -1.75%  '
$ws.Range("E33").Value = '  -0.89%  '
$ws.Range("E34").Value = '  -7.02%  '
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").Value = "'0.131"
$ws.Range("E35").Value = '  -2.70%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = "'0.119"
$ws.Range("E36").Value = '  +10.14%  '
$ws.Range("D37").Value = "'4.72"
$ws.Range("E37").Value = '  -2.01%  '
$ws.Range("E38").Value = '  -0.55%  '
$ws.Range("E39").Value = '  -3.91%  '
$ws.Range("D40").Value = "'2.71"
$ws.Range("E40").Value = '  -6.58%  '
$ws.Range("D41").Value = "'104.90"
$ws.Range("E41").Value = '  +12.44%  '
$ws.Range("E42").Value = '  -2.94%  '
$ws.Range("D43").Value = "'0.238"
$ws.Range("E43").Value = '  +1.72%  '
$ws.Range("D44").Value = "'70.98"
$ws.Range("E44").Value = '  +1.21%  '
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("D46").Value = "'12.34"
$ws.Range("E46").Value = '  +2.01%  '
$ws.Range("D47").Value = "'113.14"
$ws.Range("E47").Value = '  -5.08%  '
$ws.Range("D49").Value = "'9.12"
$ws.Range("E49").Value = '  -0.66%  '
$ws.Range("D50").Value = "'75.69"
$ws.Range("E50").Value = '  +7.07%  '
$ws.Range("E51").Value = '  +0.13%  '
